$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 59408
$ws.Range("C11").Value = "SIG-3W Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D11").Value = 388.17
$ws.Range("E11").Value = 463.78
$ws.Range("F11").Value = 26
$ws.Range("G11").Value = 10092.42
$ws.Range("B12").Value = 47438
$ws.Range("C12").Value = "SIG-3w Lilliput LED Torch &amp; Table Lamp"
$ws.Range("D12").Value = 401.81
$ws.Range("E12").Value = 480.05
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 803.62
$ws.Range("F122").Value = 6
$ws.Range("G122").Value = 314.4
$ws.Range("B126").Value = 3862.01
$ws.Range("F128").Value = 32
$ws.Range("G128").Value = 1583.36
$ws.Range("B140").Value = 58888.05
$ws.Range("B155").Value = 57756
$ws.Range("F155").Value = 60
$ws.Range("G155").Value = 3986.4
$ws.Range("B156").Value = 53925
$ws.Range("F156").Value = 1
$ws.Range("G156").Value = 66.44
$ws.Range("F167").Value = 7
$ws.Range("G167").Value = 5689.81
$ws.Range("F174").Value = 5
$ws.Range("G174").Value = 6145.75
$ws.Range("F181").Value = 36
$ws.Range("G181").Value = 3917.16
$ws.Range("B188").Value = 209541.24
$ws.Range("F223").Value = 13
$ws.Range("G223").Value = 1187.94
$ws.Range("B230").Value = 19599.21
$ws.Range("F252").Value = 235
$ws.Range("G252").Value = 4347.5
$ws.Range("B259").Value = 8197.67
$ws.Range("F275").Value = 5
$ws.Range("G275").Value = 7852.55
$ws.Range("B289").Value = 140284.22
$ws.Range("F300").Value = 40
$ws.Range("G300").Value = 11830
$ws.Range("F305").Value = 32
$ws.Range("G305").Value = 3318.4
$ws.Range("B331").Value = 214609.55
$ws.Range("F338").Value = 6
$ws.Range("G338").Value = 497.64
$ws.Range("F339").Value = 82
$ws.Range("G339").Value = 6043.4
$ws.Range("F340").Value = 132
$ws.Range("G340").Value = 10562.64
$ws.Range("F345").Value = 66
$ws.Range("G345").Value = 5422.56
$ws.Range("F346").Value = 210
$ws.Range("G346").Value = 33789
$ws.Range("F349").Value = 26
$ws.Range("G349").Value = 3820.96
$ws.Range("F350").Value = 246
$ws.Range("G350").Value = 18132.66
$ws.Range("F353").Value = 218
$ws.Range("G353").Value = 8425.700000000001
$ws.Range("F358").Value = 0
$ws.Range("G358").Value = 0
$ws.Range("F360").Value = 101
$ws.Range("G360").Value = 14497.54
$ws.Range("F363").Value = 4
$ws.Range("G363").Value = 539.12
$ws.Range("F381").Value = 206
$ws.Range("G381").Value = 22622.92
$ws.Range("F386").Value = 24
$ws.Range("G386").Value = 2982.48
$ws.Range("F387").Value = 25
$ws.Range("G387").Value = 1583.75
$ws.Range("F389").Value = 158
$ws.Range("G389").Value = 20043.88
$ws.Range("F398").Value = 2
$ws.Range("G398").Value = 438.38
$ws.Range("F401").Value = 838
$ws.Range("G401").Value = 49190.6
$ws.Range("F412").Value = 1
$ws.Range("G412").Value = 104.08
$ws.Range("F415").Value = 64
$ws.Range("G415").Value = 5318.4
$ws.Range("F417").Value = 657
$ws.Range("G417").Value = 112563.81
$ws.Range("F418").Value = 251
$ws.Range("G418").Value = 37943.67
$ws.Range("F419").Value = 9
$ws.Range("G419").Value = 3766.41
$ws.Range("F421").Value = 22
$ws.Range("G421").Value = 3527.04
$ws.Range("F423").Value = 33
$ws.Range("G423").Value = 2034.12
$ws.Range("F429").Value = 455
$ws.Range("G429").Value = 27058.85
$ws.Range("F434").Value = 141
$ws.Range("G434").Value = 20256.06
$ws.Range("B435").Value = 712652.83
$ws.Range("F437").Value = 130
$ws.Range("G437").Value = 23865.4
$ws.Range("F438").Value = 130
$ws.Range("G438").Value = 23865.4
$ws.Range("F452").Value = 32
$ws.Range("G452").Value = 7155.2
$ws.Range("B453").Value = 112314.49
$ws.Range("F484").Value = 170
$ws.Range("G484").Value = 4136.1
$ws.Range("F490").Value = 138
$ws.Range("G490").Value = 5735.28
$ws.Range("F494").Value = 108
$ws.Range("G494").Value = 709.5599999999999
$ws.Range("B507").Value = 124720.76
$ws.Range("F550").Value = 15
$ws.Range("G550").Value = 997.5
$ws.Range("F562").Value = 10
$ws.Range("G562").Value = 474
$ws.Range("F563").Value = 40
$ws.Range("G563").Value = 3684.8
$ws.Range("F564").Value = 49
$ws.Range("G564").Value = 5040.63
$ws.Range("F566").Value = 41
$ws.Range("G566").Value = 1321.84
$ws.Range("F568").Value = 48
$ws.Range("G568").Value = 3454.56
$ws.Range("F571").Value = 43
$ws.Range("G571").Value = 1237.54
$ws.Range("B577").Value = 40480.18
$ws.Range("F596").Value = 2
$ws.Range("G596").Value = 6643.5
$ws.Range("B617").Value = 246807.89
$ws.Range("F677").Value = 26
$ws.Range("G677").Value = 1389.44
$ws.Range("F681").Value = 3
$ws.Range("G681").Value = 291.51
$ws.Range("B682").Value = 1790.83
$ws.Range("F723").Value = 16
$ws.Range("G723").Value = 817.28
$ws.Range("B735").Value = 23457.05
$ws.Range("F742").Value = 2
$ws.Range("G742").Value = 31.72
$ws.Range("B743").Value = 31.72
$ws.Range("F799").Value = 119
$ws.Range("G799").Value = 10918.25
$ws.Range("B807").Value = 58526.62
$ws.Range("F843").Value = 277
$ws.Range("G843").Value = 22592.12
$ws.Range("F844").Value = 226
$ws.Range("G844").Value = 10816.36
$ws.Range("F846").Value = 93
$ws.Range("G846").Value = 14370.36
$ws.Range("F847").Value = 193
$ws.Range("G847").Value = 15741.08
$ws.Range("F848").Value = 370
$ws.Range("G848").Value = 49247
$ws.Range("F852").Value = 129
$ws.Range("G852").Value = 2801.88
$ws.Range("F863").Value = 443
$ws.Range("G863").Value = 63792
$ws.Range("F865").Value = 315
$ws.Range("G865").Value = 38023.65
$ws.Range("B867").Value = 477790.24
$ws.Range("F909").Value = 59
$ws.Range("G909").Value = 2694.53
$ws.Range("F912").Value = 1820
$ws.Range("G912").Value = 296860.2
$ws.Range("B918").Value = 330602.84
$ws.Range("B930").Value = 6061776.97
$ws.Range("B931").Value = 6061776.97

Write-Host "Applied 168 cell changes"